$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Developer name ---
$ws.Range("C3").Value = "Hudson Drozdowski"

# --- Test Case 1 (row 7): __init__ - Triangle initialized correctly ---
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "color = ""Blue""`nside_1 = 5`nside_2 = 5`nside_3 = 5"
$ws.Range("G7").Value = "Triangle initialized correctly."

# --- Test Case 2 (row 8): __init__ - Exception raised when color is blank ---
$ws.Range("E8").Value = " None"
$ws.Range("F8").Value = "color = ""   ""`nside_1 = 5`nside_2 = 5`nside_3 = 5"
$ws.Range("G8").Value = "ValueError(""Color cannot be blank."")"

# --- Test Case 3 (row 9): __init__ - Exception raised when side_1 is not numeric ---
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "color = ""Blue""`nside_1 = ""Invalid Input""`nside_2 = 5`nside_3 = 5"
$ws.Range("G9").Value = "ValueError(""Side 1 must be numeric."")"

# --- Test Case 4 (row 10): __init__ - Exception raised when side_2 is not numeric ---
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "color = ""Blue""`nside_1 = 5`nside_2 = ""Invalid Input""`nside_3 = 5"
$ws.Range("G10").Value = "ValueError(""Side 2 must be numeric."")"

# --- Test Case 5 (row 11): __init__ - Exception raised when side_3 is not numeric ---
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "color = ""Blue""`nside_1 = 5`nside_2 = 5`nside_3 = ""Invalid Input"""
$ws.Range("G11").Value = "ValueError(""Side 3 must be numeric."")"

# --- Test Case 6 (row 12): __str__ - Returns string formatted appropriately ---
$ws.Range("E12").Value = "color = ""Blue""`nside_1 = 5`nside_2 = 5`nside_3 = 5"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "The shape color is Blue.`nThis triangle has three sides with the lengths of 5, 5, 5 centimeters."

# --- Test Case 7 (row 13): calculate_area - Returns correct calculated value ---
$ws.Range("E13").Value = "color = ""Blue""`nside_1 = 5`nside_2 = 5`nside_3 = 5"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = 10.83

# --- Test Case 8 (row 14): calculate_perimeter - Returns correct calculated value ---
$ws.Range("E14").Value = "color = ""Blue""`nside_1 = 5`nside_2 = 5`nside_3 = 5"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = 15
